# Análise de Valor Agregado - atualização do burndown e valor agregado
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Preenche as horas trabalhadas (coluna B) das tarefas que foram concluídas
$ws.Range("B43").Value = 8
$ws.Range("B44").Value = 4
$ws.Range("B47").Value = 3.5
$ws.Range("B49").Value = 0.5

# Marca essas tarefas como concluídas na Sprint 6 (coluna I)
$ws.Range("I43").Value = "S"
$ws.Range("I44").Value = "S"
$ws.Range("I47").Value = "S"
$ws.Range("I49").Value = "S"

# Remove as planilhas vazias Plan2 e Plan3
$wb.Worksheets.Item("Plan2").Delete()
$wb.Worksheets.Item("Plan3").Delete()

# Seleciona a célula B43 na visão da planilha ativa
$ws.Range("B43").Select()

$wb.Save()
